# Apply MOSIP "individual_type" MEC prod data update.
# The data table is rebuilt with a reshuffled column order (code, name,
# lang_code, is_active) plus six new audit columns (cr_by, cr_dtimes,
# upd_by, upd_dtimes, is_deleted, del_dtimes), and the header font is no
# longer bold.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row (row 1): new column order + six new trailing audit columns.
# Un-bold the header font (still centered, bordered, D1 kept as text fmt).
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "lang_code"
$ws.Range("D1").Value = "is_active"
$ws.Range("E1").Value = "cr_by"
$ws.Range("F1").Value = "cr_dtimes"
$ws.Range("G1").Value = "upd_by"
$ws.Range("H1").Value = "upd_dtimes"
$ws.Range("I1").Value = "is_deleted"
$ws.Range("J1").Value = "del_dtimes"

$ws.Range("A1:D1").Font.Bold = $false

# ---------------------------------------------------------------------
# Data rows 2-7: code / name / lang_code in the new column order.
# ---------------------------------------------------------------------
$rows = @(
    @{ Code = "FR";  Name = "Foreigner";      Lang = "eng" },
    @{ Code = "NFR"; Name = "Non-Foreigner";  Lang = "eng" },
    @{ Code = "FR";  Name = "Étranger";       Lang = "fra" },
    @{ Code = "NFR"; Name = "Non-étranger";   Lang = "fra" },
    @{ Code = "FR";  Name = "أجنبي";          Lang = "ara" },
    @{ Code = "NFR"; Name = "غير أجنبي";      Lang = "ara" }
)

$creationDate = 44776.353977430554

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $ws.Range("A$r").Value = $row.Code
    $ws.Range("B$r").Value = $row.Name
    $ws.Range("C$r").Value = $row.Lang

    if ($r -le 5) {
        # eng/fra rows: is_active stored as the literal text "true"
        $ws.Range("D$r").NumberFormat = "@"
        $ws.Range("D$r").Value = "'true"
    } else {
        # ara rows: is_active stored as a real boolean TRUE
        $ws.Range("D$r").Value = $true
    }

    $ws.Range("E$r").Value = "rediet"

    $ws.Range("F$r").Value = $creationDate
    $ws.Range("F$r").NumberFormat = "mm:ss.0"

    $ws.Range("I$r").Value = $false
}

$null = $ws.Range("H12").Select()

Write-Host "individual_type sheet updated with MEC prod data"
